# Update the cryptocurrency price/volume table on Sheet1 with the latest
# scraped values (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper note: many "Price" values in column D look numeric (e.g. "603.87")
# but must remain plain text, exactly as authored (matching the workbook's
# existing inlineStr cells). Assigning them directly as ".Value" causes
# Excel to auto-convert them into real numbers (losing trailing zeros /
# exact formatting), so we prefix with a single quote to force text entry
# and then reset the cell style back to Normal so no quote-prefix /
# "@" text-format style lingers on the cell.

$ws.Range("D2").Value = "'69.130.61"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.11%  "

$ws.Range("D3").Value = "'2.742.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.44%  "

$ws.Range("D5").Value = "'603.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.07%  "

$ws.Range("D6").Value = "'165.83"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.63%  "

$ws.Range("D8").Value = "'0.547"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.64%  "

$ws.Range("D9").Value = "'2.741.06"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.43%  "

$ws.Range("E10").Value = "  +0.08%  "

$ws.Range("E11").Value = "  +3.53%  "

$ws.Range("E12").Value = "  -0.06%  "

$ws.Range("E13").Value = "  +1.26%  "

$ws.Range("D14").Value = "'28.72"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.02%  "

$ws.Range("D15").Value = "'3.241.83"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.58%  "

$ws.Range("D16").Value = "'0.0000191"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.47%  "

$ws.Range("D17").Value = "'68.980.16"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.14%  "

$ws.Range("D18").Value = "'2.741.46"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.51%  "

$ws.Range("E19").Value = "  +4.64%  "

$ws.Range("E20").Value = "  +5.34%  "

$ws.Range("D21").Value = "'368.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.11%  "

$ws.Range("D22").Value = "'4.57"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.09%  "

$ws.Range("D23").Value = "'4.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.48%  "

$ws.Range("E24").Value = "  +3.63%  "

$ws.Range("D25").Value = "'74.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.59%  "

$ws.Range("E26").Value = "  -0.09%  "

$ws.Range("D27").Value = "'10.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.89%  "

$ws.Range("E28").Value = "  +2.48%  "

$ws.Range("E29").Value = "  +2.22%  "

$ws.Range("D30").Value = "'603.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.91%  "

$ws.Range("D31").Value = "'1.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.09%  "

$ws.Range("E32").Value = "  +3.89%  "

$ws.Range("E33").Value = "  +3.45%  "

$ws.Range("D34").Value = "'1.97"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.72%  "

$ws.Range("E35").Value = "  +3.51%  "

$ws.Range("E36").Value = "  +4.75%  "

$ws.Range("E37").Value = "  +0.09%  "

$ws.Range("D38").Value = "'163.11"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.15%  "

$ws.Range("D39").Value = "'20.12"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.37%  "

$ws.Range("E40").Value = "  +3.11%  "

$ws.Range("D41").Value = "'1.93"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.47%  "

$ws.Range("D42").Value = "'5.50"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.89%  "

$ws.Range("E43").Value = "  +3.70%  "

$ws.Range("E44").Value = "  +1.35%  "

$ws.Range("E45").Value = "  -4.35%  "

$ws.Range("E46").Value = "  +0.06%  "

$ws.Range("D47").Value = "'159.44"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.43%  "

$ws.Range("E48").Value = "  +5.61%  "

$ws.Range("E49").Value = "  +6.65%  "

$ws.Range("E50").Value = "  +7.90%  "

$ws.Range("D51").Value = "'22.25"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.28%  "
